$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (for the week of 2022-04-14, serial 44665) was
# inserted as row 29. Every existing record previously on rows 29-74 shifts
# down by one row (to rows 30-75); the insert preserves the D-column date
# style carried by the rest of the table.
$ws.Rows.Item(29).Insert()

$ws.Cells.Item(29, 1).Value = 1
$ws.Cells.Item(29, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(29, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(29, 4).Value = 44665
$ws.Cells.Item(29, 5).Value = 15
$ws.Cells.Item(29, 6).Value = 100112040
$ws.Cells.Item(29, 7).Value = "Cilantro"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 300
$ws.Cells.Item(29, 11).Value = 1900
$ws.Cells.Item(29, 12).Value = 2000
$ws.Cells.Item(29, 13).Value = 1950
$ws.Cells.Item(29, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(29, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(29, 16).Value = 975
$ws.Cells.Item(29, 17).Value = 2
$ws.Cells.Item(29, 18).Value = "Hortaliza"
